# Generate Report for Handback
# c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b has now been handed back (in sync with en-US).
# 17d16921-d803-4efd-bb22-ef1ca06a2a3e is still ready for handoff (unchanged status).
# This script rewrites the three worksheets (Overview, zh-cn, de-de) to reflect the
# new handback status, swapping row 2 / row 3 ordering (c2ffc5e8 now first) and
# populating the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for c2ffc5e8 on the language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hyperlink target URLs (reused from the file's existing relationships so the
# same files keep resolving to the same GitHub blobs).
# ---------------------------------------------------------------------------
$url_17d_md       = "https://github.com/OpenLocalizationTest/oltest/blob/c463897b4bc5670b746f73360af9d2c57a535b03/e2e/17d16921-d803-4efd-bb22-ef1ca06a2a3e.md"
$url_c2f_md       = "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md"
$url_config       = "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/.localization-config"

$url_17d_zhcn_xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8341175a15e72d81211d1263b4bae898989241a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.zh-cn.xlf"
$url_c2f_zhcn_xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/26d8541d6230e3f3e5ecaf821cd8448033cfcc1a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf"

$url_17d_dede_xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d94394ebae535fe0b9027c2b4b735700a647dc7a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.de-de.xlf"
$url_c2f_dede_xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95256d14e2226de8f26c490d832cdc5efc7bc2e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf"

$name_17d_md  = "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md"
$name_c2f_md  = "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md"
$name_config  = ".localization-config"

$name_17d_zhcn_xlf = "17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.zh-cn.xlf"
$name_c2f_zhcn_xlf = "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf"
$name_17d_dede_xlf = "17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.de-de.xlf"
$name_c2f_dede_xlf = "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf"

$status_handed_back = "Handed back: in sync with en-US"
$status_ready       = "Ready for handoff"
$status_not_loc     = "Not to be localized"
$status_include     = "Include"
$status_ignored     = "Ignored"
$dt_epoch           = "0001-01-01 00:00:00"

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Clear existing hyperlinks on A2:A4 so we can recreate them in the right order.
$wsOverview.Range("A2:A4").Hyperlinks.Delete()

# Row 2 is now c2ffc5e8 (handed back).
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $url_c2f_md, "", "", $name_c2f_md)
$wsOverview.Range("B2").Value = $status_handed_back
$wsOverview.Range("C2").Value = $status_handed_back

# Row 3 is now 17d16921 (still ready for handoff).
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $url_17d_md, "", "", $name_17d_md)
$wsOverview.Range("B3").Value = $status_ready
$wsOverview.Range("C3").Value = $status_ready

# Row 4 (.localization-config) is unchanged, but re-add so ordering / rIds stay sane.
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $url_config, "", "", $name_config)
$wsOverview.Range("B4").Value = $status_not_loc
$wsOverview.Range("C4").Value = $status_not_loc

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2:A4").Hyperlinks.Delete()
$wsZh.Range("C2:C4").Hyperlinks.Delete()
$wsZh.Range("E2:E4").Hyperlinks.Delete()
$wsZh.Range("F2:F4").Hyperlinks.Delete()

# Row 2: c2ffc5e8 (handed back) zh-cn details.
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $url_c2f_md, "", "", $name_c2f_md)
$wsZh.Range("B2").Value = $status_handed_back
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $url_c2f_zhcn_xlf, "", "", $name_c2f_zhcn_xlf)
$wsZh.Range("D2").Value = "2016-03-10 03:10:20"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $url_c2f_md, "", "", $name_c2f_md)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $url_c2f_zhcn_xlf, "", "", $name_c2f_zhcn_xlf)
$wsZh.Range("G2").Value = "2016-03-10 03:11:09"
$wsZh.Range("H2").Value = $status_include

# Row 3: 17d16921 (still ready for handoff) zh-cn details.
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $url_17d_md, "", "", $name_17d_md)
$wsZh.Range("B3").Value = $status_ready
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $url_17d_zhcn_xlf, "", "", $name_17d_zhcn_xlf)
$wsZh.Range("D3").Value = "2016-03-10 03:09:36"
$wsZh.Range("G3").Value = $dt_epoch
$wsZh.Range("H3").Value = $status_include

# Row 4: .localization-config (unchanged).
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $url_config, "", "", $name_config)
$wsZh.Range("B4").Value = $status_not_loc
$wsZh.Range("D4").Value = $dt_epoch
$wsZh.Range("G4").Value = $dt_epoch
$wsZh.Range("H4").Value = $status_ignored

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2:A4").Hyperlinks.Delete()
$wsDe.Range("C2:C4").Hyperlinks.Delete()
$wsDe.Range("E2:E4").Hyperlinks.Delete()
$wsDe.Range("F2:F4").Hyperlinks.Delete()

# Row 2: c2ffc5e8 (handed back) de-de details.
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $url_c2f_md, "", "", $name_c2f_md)
$wsDe.Range("B2").Value = $status_handed_back
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $url_c2f_dede_xlf, "", "", $name_c2f_dede_xlf)
$wsDe.Range("D2").Value = "2016-03-10 03:10:23"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $url_c2f_md, "", "", $name_c2f_md)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $url_c2f_dede_xlf, "", "", $name_c2f_dede_xlf)
$wsDe.Range("G2").Value = "2016-03-10 03:11:14"
$wsDe.Range("H2").Value = $status_include

# Row 3: 17d16921 (still ready for handoff) de-de details.
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $url_17d_md, "", "", $name_17d_md)
$wsDe.Range("B3").Value = $status_ready
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $url_17d_dede_xlf, "", "", $name_17d_dede_xlf)
$wsDe.Range("D3").Value = "2016-03-10 03:09:39"
$wsDe.Range("G3").Value = $dt_epoch
$wsDe.Range("H3").Value = $status_include

# Row 4: .localization-config (unchanged).
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $url_config, "", "", $name_config)
$wsDe.Range("B4").Value = $status_not_loc
$wsDe.Range("D4").Value = $dt_epoch
$wsDe.Range("G4").Value = $dt_epoch
$wsDe.Range("H4").Value = $status_ignored

Write-Host "Handback report regenerated."
